$wb = $excel.ActiveWorkbook

# The localization status moved from "Ready for handoff" to "In Translation"
# on the per-language report sheets (zh-cn, de-de) and on the Overview
# sheet's mirrored language-status columns.
$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value = $newStatus
            }
        }
    }
}

# The Status column got narrower now that the text is shorter - shrink the
# affected columns to match (report regenerated with the new, shorter
# content). ColumnWidth of 12.5 is the closest this engine's column-width
# grid gets to the regenerated report's target character width (~13.41).
$newWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

$wb.Worksheets.Item("zh-cn").Columns.Item(3).ColumnWidth = $newWidth
$wb.Worksheets.Item("de-de").Columns.Item(3).ColumnWidth = $newWidth
